# Edit script for khl_stats_1369_ext.xlsx
# Commit: chore(runtime): publish files + archive (2025-11-10 11:06:44)
#
# 1) Matches_SOG: append 3 new match rows (452-454) for 2025-11-09 games.
# 2) Shots_HA: refresh as_of_utc snapshot timestamp and the per-team stat
#    values that shifted because of the newly completed games.
# 3) Shots_Summary: same refresh as Shots_HA, on the summary sheet.
# 4) Meta_ext: bump as_of_utc and build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matches_SOG - append new match rows
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @(452, "897749", "2025-11-09T10:00:00", "Амур",     "Автомобилист",  33, 39),
    @(453, "897751", "2025-11-09T10:00:00", "Адмирал",  "Трактор",       36, 26),
    @(454, "897750", "2025-11-09T15:00:00", "Барыс",    "Металлург Мг",  32, 34)
)

foreach ($row in $newMatches) {
    $r = $row[0]
    $wsMatches.Cells.Item($r, 1).Value = "'" + $row[1]
    $wsMatches.Cells.Item($r, 2).Value = $row[2]
    $wsMatches.Cells.Item($r, 3).Value = $row[3]
    $wsMatches.Cells.Item($r, 4).Value = $row[4]
    $wsMatches.Cells.Item($r, 5).Value = $row[5]
    $wsMatches.Cells.Item($r, 6).Value = $row[6]
    $wsMatches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# 2) Shots_HA - refresh as_of_utc + derived shot stats
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_HA")

$ws.Cells.Item(2, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(3, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(3, 6).Value = 26
$ws.Cells.Item(3, 11).Value = 731
$ws.Cells.Item(3, 12).Value = 798
$ws.Cells.Item(3, 13).Value = 28.1
$ws.Cells.Item(3, 14).Value = 30.7
$ws.Cells.Item(4, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(4, 5).Value = 17
$ws.Cells.Item(4, 7).Value = 653
$ws.Cells.Item(4, 8).Value = 462
$ws.Cells.Item(4, 9).Value = 38.4
$ws.Cells.Item(5, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(6, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(6, 5).Value = 20
$ws.Cells.Item(6, 7).Value = 614
$ws.Cells.Item(6, 8).Value = 710
$ws.Cells.Item(6, 9).Value = 30.7
$ws.Cells.Item(6, 10).Value = 35.5
$ws.Cells.Item(7, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(7, 5).Value = 30
$ws.Cells.Item(7, 7).Value = 957
$ws.Cells.Item(7, 8).Value = 953
$ws.Cells.Item(7, 10).Value = 31.8
$ws.Cells.Item(8, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(9, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(10, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(11, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(12, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(13, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(13, 6).Value = 17
$ws.Cells.Item(13, 11).Value = 488
$ws.Cells.Item(13, 12).Value = 460
$ws.Cells.Item(13, 13).Value = 28.7
$ws.Cells.Item(13, 14).Value = 27.1
$ws.Cells.Item(14, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(15, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(16, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(17, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(18, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(19, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(20, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(21, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(21, 6).Value = 25
$ws.Cells.Item(21, 11).Value = 852
$ws.Cells.Item(21, 12).Value = 820
$ws.Cells.Item(21, 13).Value = 34.1
$ws.Cells.Item(21, 14).Value = 32.8
$ws.Cells.Item(22, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(23, 4).Value = '2025-11-09T15:00:00Z'

# ---------------------------------------------------------------------
# 3) Shots_Summary - refresh as_of_utc + derived shot stats
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_Summary")

$ws.Cells.Item(2, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(3, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(3, 5).Value = 44
$ws.Cells.Item(3, 6).Value = 1252
$ws.Cells.Item(3, 7).Value = 1356
$ws.Cells.Item(3, 8).Value = 28.5
$ws.Cells.Item(4, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(4, 5).Value = 37
$ws.Cells.Item(4, 6).Value = 1291
$ws.Cells.Item(4, 7).Value = 1023
$ws.Cells.Item(4, 9).Value = 27.6
$ws.Cells.Item(5, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(6, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(6, 5).Value = 41
$ws.Cells.Item(6, 6).Value = 1206
$ws.Cells.Item(6, 7).Value = 1479
$ws.Cells.Item(6, 8).Value = 29.4
$ws.Cells.Item(6, 9).Value = 36.1
$ws.Cells.Item(7, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(7, 5).Value = 45
$ws.Cells.Item(7, 6).Value = 1372
$ws.Cells.Item(7, 7).Value = 1468
$ws.Cells.Item(8, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(9, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(10, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(11, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(12, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(13, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(13, 5).Value = 42
$ws.Cells.Item(13, 6).Value = 1372
$ws.Cells.Item(13, 7).Value = 1106
$ws.Cells.Item(13, 8).Value = 32.7
$ws.Cells.Item(13, 9).Value = 26.3
$ws.Cells.Item(14, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(15, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(16, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(17, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(18, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(19, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(20, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(21, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(21, 5).Value = 43
$ws.Cells.Item(21, 6).Value = 1451
$ws.Cells.Item(21, 7).Value = 1361
$ws.Cells.Item(21, 8).Value = 33.7
$ws.Cells.Item(21, 9).Value = 31.7
$ws.Cells.Item(22, 4).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(23, 4).Value = '2025-11-09T15:00:00Z'

# ---------------------------------------------------------------------
# 4) Meta_ext - bump as_of_utc + build_version
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Meta_ext")

$ws.Cells.Item(2, 2).Value = '2025-11-09T15:00:00Z'
$ws.Cells.Item(2, 4).Value = 57
